$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data (price & 1h volume change); rows 42-51 also have
# coin name/link shifts because "Frax" dropped out of the top list and
# "BabyDogeCoin" entered at the bottom.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.636.82'
$ws.Range('E2').Value = '  +1.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.868.72'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.43%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '331.60'
$ws.Range('E5').Value = '  +2.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4723'
$ws.Range('E7').Value = '  +4.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3938'
$ws.Range('E8').Value = '  +1.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.96'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08045'
$ws.Range('E10').Value = '  +1.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.024'
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.00'
$ws.Range('E12').Value = '  +2.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.884.36'
$ws.Range('E13').Value = '  +1.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.947'
$ws.Range('E14').Value = '  +0.70%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.135'
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.007'
$ws.Range('E16').Value = '  +0.69%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001046'
$ws.Range('E17').Value = '  +1.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '86.76'
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06652'
$ws.Range('E19').Value = '  +2.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.19'
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '27.653.34'
$ws.Range('E22').Value = '  +1.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.501'
$ws.Range('E23').Value = '  -0.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.96'
$ws.Range('E24').Value = '  +0.81%  '
$ws.Range('E25').Value = '  +1.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.091.24'
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '158.85'
$ws.Range('E27').Value = '  +3.89%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.24'
$ws.Range('E28').Value = '  +2.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.095'
$ws.Range('E29').Value = '  +1.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.551'
$ws.Range('E30').Value = '  +0.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '122.24'
$ws.Range('E31').Value = '  +1.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9688'
$ws.Range('E32').Value = '  +3.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09521'
$ws.Range('E33').Value = '  +2.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.448'
$ws.Range('E34').Value = '  -3.42%  '
$ws.Range('E35').Value = '  -0.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.327'
$ws.Range('E36').Value = '  +0.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06100'
$ws.Range('E37').Value = '  +1.60%  '
$ws.Range('E38').Value = '  +0.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.230'
$ws.Range('E39').Value = '  +0.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.138'
$ws.Range('E40').Value = '  -1.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6019'
$ws.Range('E41').Value = '  +1.74%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1899'
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.25'
$ws.Range('E43').Value = '  +1.03%  '
$ws.Range('B44').Value = 'WEMIXTOKEN'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.256'
$ws.Range('E44').Value = '  -1.43%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5709'
$ws.Range('E45').Value = '  +1.10%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.22'
$ws.Range('E46').Value = '  +1.29%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.940'
$ws.Range('E47').Value = '  +0.54%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.383'
$ws.Range('E48').Value = '  +0.37%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06865'
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '114.40'
$ws.Range('E50').Value = '  +5.82%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.00000000299'
$ws.Range('E51').Value = '  +12.99%  '
